# Appends rows 31-38 (the new bug-report data) to Sheet1, extending the
# used range from A1:T30 to A1:T38 as in the target workbook.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31
$ws.Range("A31").Value = 29
$ws.Range("B31").Value = 20
$ws.Range("C31").Value = 0
$ws.Range("D31").Value = 0.003
$ws.Range("E31").Value = "FGSM"
$ws.Range("F31").Value = ""
$ws.Range("G31").Value = ""
$ws.Range("H31").Value = ""
$ws.Range("I31").Value = "<function relu at 0x128d789d8>"
$ws.Range("J31").Value = 0.9329000115394592
$ws.Range("K31").Value = 0.7416999936103821
$ws.Range("L31").Value = 0.5971999764442444
$ws.Range("M31").Value = 0.2206476628780365
$ws.Range("N31").Value = 0.74456787109375
$ws.Range("O31").Value = 0.7416999936103821
$ws.Range("P31").Value = "logs/results_328.log"
$ws.Range("Q31").Value = "weights/model_328.ckpt"
$ws.Range("R31").Value = "tb/328"
$ws.Range("S31").Value = "(12.799518, 16.56418, 25.464622, 25.766693, 14.9766655, 6.904501, 3.7851212)"
$ws.Range("T31").Value = "(309.03058, 13.829748, 16.421766, 12.924315, 13.326067, 11.006118, 10.691159, 13.276878)"

# Row 32
$ws.Range("A32").Value = 30
$ws.Range("B32").Value = 30
$ws.Range("C32").Value = 0.0005
$ws.Range("D32").Value = 0.003
$ws.Range("E32").Value = "Regular"
$ws.Range("F32").Value = ""
$ws.Range("G32").Value = ""
$ws.Range("H32").Value = ""
$ws.Range("I32").Value = "<function relu at 0x113fe39d8>"
$ws.Range("J32").Value = 0.9472000002861023
$ws.Range("K32").Value = 0.3325000107288361
$ws.Range("L32").Value = 0.1923000067472458
$ws.Range("M32").Value = 0.1799724549055099
$ws.Range("N32").Value = 3.20237398147583
$ws.Range("O32").Value = 0.3325000107288361
$ws.Range("P32").Value = "logs/results_331.log"
$ws.Range("Q32").Value = "weights/model_331.ckpt"
$ws.Range("R32").Value = "tb/331"
$ws.Range("S32").Value = "(1.8040013, 1.8094578, 1.9554862, 2.5528448, 3.399231, 4.635884, 6.39143)"
$ws.Range("T32").Value = "(28.089184, 4.674623, 3.6114936, 3.2815678, 3.1776936, 2.116904, 2.2179012, 2.6026235)"

# Row 33
$ws.Range("A33").Value = 31
$ws.Range("B33").Value = 30
$ws.Range("C33").Value = 0.0009
$ws.Range("D33").Value = 0.003
$ws.Range("E33").Value = "Regular"
$ws.Range("F33").Value = ""
$ws.Range("G33").Value = ""
$ws.Range("H33").Value = ""
$ws.Range("I33").Value = "<function relu at 0x118f1b9d8>"
$ws.Range("J33").Value = 0.8956000208854675
$ws.Range("K33").Value = 0.3217999935150146
$ws.Range("L33").Value = 0.2044000029563904
$ws.Range("M33").Value = 0.3738516569137573
$ws.Range("N33").Value = 3.476548433303833
$ws.Range("O33").Value = 0.3217999935150146
$ws.Range("P33").Value = "logs/results_332.log"
$ws.Range("Q33").Value = "weights/model_332.ckpt"
$ws.Range("R33").Value = "tb/332"
$ws.Range("S33").Value = "(1.4544966, 1.1530415, 1.2943891, 2.1282125, 4.040144, 7.1142826, 6.1251535)"
$ws.Range("T33").Value = "(21.872171, 4.500764, 3.899859, 2.0486226, 2.4281611, 2.2534535, 1.7587844, 2.4474022)"

# Row 34
$ws.Range("A34").Value = 32
$ws.Range("B34").Value = 60
$ws.Range("C34").Value = 0.0015
$ws.Range("D34").Value = 0.0003
$ws.Range("E34").Value = "Regular"
$ws.Range("F34").Value = ""
$ws.Range("G34").Value = ""
$ws.Range("H34").Value = ""
$ws.Range("I34").Value = "<function relu at 0x10b0f09d8>"
$ws.Range("J34").Value = 0.9347000122070312
$ws.Range("K34").Value = 0.2572999894618988
$ws.Range("L34").Value = 0.1509999930858612
$ws.Range("M34").Value = 0.2364677786827087
$ws.Range("N34").Value = 4.247204780578613
$ws.Range("O34").Value = 0.2572999894618988
$ws.Range("P34").Value = "logs/results_333.log"
$ws.Range("Q34").Value = "weights/model_333.ckpt"
$ws.Range("R34").Value = "tb/333"
$ws.Range("S34").Value = "(0.4250503, 0.4371529, 0.74407357, 1.2906153, 2.5302472, 4.9194384, 6.6656017)"
$ws.Range("T34").Value = "(6.92537, 5.7589846, 3.67073, 3.039194, 3.1920185, 5.1396704, 3.4247825, 4.5205474)"

# Row 35
$ws.Range("A35").Value = 33
$ws.Range("B35").Value = 60
$ws.Range("C35").Value = 0.000015
$ws.Range("D35").Value = 0.0003
$ws.Range("E35").Value = "Regular"
$ws.Range("F35").Value = ""
$ws.Range("G35").Value = ""
$ws.Range("H35").Value = ""
$ws.Range("I35").Value = "<function relu at 0x1113b09d8>"
$ws.Range("J35").Value = 0.9394999742507935
$ws.Range("K35").Value = 0.09109999984502792
$ws.Range("L35").Value = 0.03319999948143959
$ws.Range("M35").Value = 0.221238300204277
$ws.Range("N35").Value = 6.377280712127686
$ws.Range("O35").Value = 0.09109999984502792
$ws.Range("P35").Value = "logs/results_336.log"
$ws.Range("Q35").Value = "weights/model_336.ckpt"
$ws.Range("R35").Value = "tb/336"
$ws.Range("S35").Value = "(2.0724185, 1.3095005, 1.6185311, 1.3509899, 1.5398465, 1.9591715, 4.275175)"
$ws.Range("T35").Value = "(39.100674, 3.319005, 3.0519383, 1.6382135, 2.70516, 2.7627785, 5.8854613, 6.491546)"

# Row 36
$ws.Range("A36").Value = 34
$ws.Range("B36").Value = 30
$ws.Range("C36").Value = 0.000023
$ws.Range("D36").Value = 0.0003
$ws.Range("E36").Value = "Regular"
$ws.Range("F36").Value = ""
$ws.Range("G36").Value = ""
$ws.Range("H36").Value = ""
$ws.Range("I36").Value = "<function relu at 0x1132949d8>"
$ws.Range("J36").Value = 0.9189000129699707
$ws.Range("K36").Value = 0.05869999900460243
$ws.Range("L36").Value = 0.02209999971091747
$ws.Range("M36").Value = 0.3333184719085693
$ws.Range("N36").Value = 6.367059707641602
$ws.Range("O36").Value = 0.05869999900460243
$ws.Range("P36").Value = "logs/results_346.log"
$ws.Range("Q36").Value = "weights/model_346.ckpt"
$ws.Range("R36").Value = "tb/346"
$ws.Range("S36").Value = "(1.5257524, 0.7150309, 0.56886697, 0.49205807, 0.57274866, 1.5198021, 4.5895967)"
$ws.Range("T36").Value = "(33.718807, 2.9477801, 2.0959005, 1.7632871, 2.2024982, 4.748903, 8.546004, 6.910356)"

# Row 37
$ws.Range("A37").Value = 35
$ws.Range("B37").Value = 100
$ws.Range("C37").Value = 0.007
$ws.Range("D37").Value = 0.0003
$ws.Range("E37").Value = "Regular"
$ws.Range("F37").Value = ""
$ws.Range("G37").Value = ""
$ws.Range("H37").Value = ""
$ws.Range("I37").Value = "<function relu at 0x1226c79d8>"
$ws.Range("J37").Value = 0.9646999835968018
$ws.Range("K37").Value = 0.1290999948978424
$ws.Range("L37").Value = 0.01889999955892563
$ws.Range("M37").Value = 0.1261469274759293
$ws.Range("N37").Value = 6.858481407165527
$ws.Range("O37").Value = 0.1290999948978424
$ws.Range("P37").Value = "logs/results_353.log"
$ws.Range("Q37").Value = "weights/model_353.ckpt"
$ws.Range("R37").Value = "tb/353"
$ws.Range("S37").Value = "(0.6469947, 0.73037505, 0.9719021, 1.3849623, 2.158197, 3.6967816, 6.6338224)"
$ws.Range("T37").Value = "(9.72191, 3.5845523, 3.115291, 2.8329852, 2.7452235, 2.6821659, 2.723758, 2.6171021)"

# Row 38
$ws.Range("A38").Value = 36
$ws.Range("B38").Value = 100
$ws.Range("C38").Value = 0.07
$ws.Range("D38").Value = 0.0003
$ws.Range("E38").Value = "Regular"
$ws.Range("F38").Value = ""
$ws.Range("G38").Value = ""
$ws.Range("H38").Value = ""
$ws.Range("I38").Value = "<function relu at 0x1226499d8>"
$ws.Range("J38").Value = 0.9384999871253967
$ws.Range("K38").Value = 0.3431999981403351
$ws.Range("L38").Value = 0.217399999499321
$ws.Range("M38").Value = 0.2118040025234222
$ws.Range("N38").Value = 2.253601551055908
$ws.Range("O38").Value = 0.3431999981403351
$ws.Range("P38").Value = "logs/results_354.log"
$ws.Range("Q38").Value = "weights/model_354.ckpt"
$ws.Range("R38").Value = "tb/354"
$ws.Range("S38").Value = "(0.2376155, 0.2859654, 0.43169716, 0.67097735, 1.0554945, 1.7692584, 2.8368058)"
$ws.Range("T38").Value = "(3.0683854, 1.7650508, 1.7416625, 1.7273158, 1.7367424, 1.746138, 1.7387878, 1.7745363)"
